$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.ClearFormats()
}

# Price (column D) updates - forced to text to preserve original formatting
Set-TextValue "D2" "25.904.09"
Set-TextValue "D3" "1.740.79"
Set-TextValue "D5" "248.11"
Set-TextValue "D6" "1.000"
Set-TextValue "D7" "0.5114"
Set-TextValue "D8" "0.2741"
Set-TextValue "D9" "0.06183"
Set-TextValue "D10" "1.738.95"
Set-TextValue "D11" "0.07231"
Set-TextValue "D12" "15.12"
Set-TextValue "D13" "0.6475"
Set-TextValue "D15" "77.60"
Set-TextValue "D16" "1.000"
Set-TextValue "D18" "25.923.77"
Set-TextValue "D19" "11.81"
Set-TextValue "D20" "0.000006806"
Set-TextValue "D21" "1.962.45"
Set-TextValue "D22" "4.275"
Set-TextValue "D23" "8.642"
Set-TextValue "D24" "5.389"
Set-TextValue "D25" "136.45"
Set-TextValue "D26" "1.500"
Set-TextValue "D27" "15.22"
Set-TextValue "D28" "1.773"
Set-TextValue "D29" "105.39"
Set-TextValue "D30" "3.910"
Set-TextValue "D31" "0.08220"
Set-TextValue "D32" "3.643"
Set-TextValue "D33" "0.04679"
Set-TextValue "D34" "2.656"
Set-TextValue "D35" "0.9983"
Set-TextValue "D36" "0.6259"
Set-TextValue "D37" "2.728"
Set-TextValue "D38" "0.01602"
Set-TextValue "D39" "1.921"
Set-TextValue "D40" "0.9999"
Set-TextValue "D41" "99.99"
Set-TextValue "D44" "4.989"
Set-TextValue "D46" "6.287"
Set-TextValue "D47" "55.30"
Set-TextValue "D48" "0.05236"
Set-TextValue "D49" "30.64"
Set-TextValue "D50" "7.516"
Set-TextValue "D51" "0.3411"

# Coin name / Link / Volume updates
$ws.Range("E2").Value = "  -0.67%  "
$ws.Range("E3").Value = "  -0.50%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("E5").Value = "  +5.59%  "
$ws.Range("E6").Value = "  +0.12%  "
$ws.Range("E7").Value = "  -3.51%  "
$ws.Range("E8").Value = "  -1.78%  "
$ws.Range("E9").Value = "  -0.03%  "
$ws.Range("E10").Value = "  -0.40%  "
$ws.Range("E11").Value = "  +0.75%  "
$ws.Range("E12").Value = "  -1.79%  "
$ws.Range("E13").Value = "  +0.42%  "
$ws.Range("E14").Value = "  +0.07%  "
$ws.Range("E15").Value = "  -1.02%  "
$ws.Range("E16").Value = "  +0.13%  "
$ws.Range("E17").Value = "  +0.15%  "
$ws.Range("E18").Value = "  -0.25%  "
$ws.Range("E19").Value = "  +1.10%  "
$ws.Range("E20").Value = "  +1.25%  "
$ws.Range("E21").Value = "  -0.35%  "
$ws.Range("E22").Value = "  -0.83%  "
$ws.Range("E23").Value = "  -1.25%  "
$ws.Range("E24").Value = "  +3.11%  "
$ws.Range("E25").Value = "  -1.22%  "
$ws.Range("E26").Value = "  -0.76%  "
$ws.Range("E27").Value = "  -0.48%  "
$ws.Range("E28").Value = "  -1.80%  "
$ws.Range("E29").Value = "  +0.80%  "
$ws.Range("E30").Value = "  +2.86%  "
$ws.Range("E31").Value = "  -0.70%  "
$ws.Range("E32").Value = "  -0.44%  "
$ws.Range("E33").Value = "  +2.14%  "
$ws.Range("E34").Value = "  +0.67%  "
$ws.Range("E35").Value = "  -0.58%  "
$ws.Range("E36").Value = "  -1.04%  "
$ws.Range("E37").Value = "  +0.46%  "
$ws.Range("E39").Value = "  -1.32%  "
$ws.Range("E40").Value = "  +0.14%  "
$ws.Range("E41").Value = "  +0.17%  "
$ws.Range("E42").Value = "  +1.59%  "
$ws.Range("E43").Value = "  -2.04%  "
$ws.Range("E44").Value = "  -0.79%  "
$ws.Range("E45").Value = "  -1.51%  "
$ws.Range("E46").Value = "  -0.62%  "
$ws.Range("E47").Value = "  +2.10%  "
$ws.Range("E48").Value = "  -2.14%  "
$ws.Range("E49").Value = "  -0.64%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("E50").Value = "  -1.41%  "
$ws.Range("B51").Value = "Decentraland"
$ws.Range("C51").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("E51").Value = "  -1.39%  "
